# Config.xlsx edit: add OCR / Document Understanding configuration rows
# and a new "Description" column explaining each setting.
#
# Original rows (A=Name, B=Value):
#   1 Name | Value
#   2 DirectoryPath
#   3 EndProcessMessage
#   4 ApiKey
#   5 Endpoint
#
# Target rows (A=Name, B=Value, C=Description):
#   1 Name | Value | Description
#   2 DirectoryPath            | The path of the directory in which the files to be processed are stored
#   3 EndProcessMessage        | The message that should be displayed when the workflow execution ends
#   4 OCRApiKey                | The OCR Engine API Key
#   5 DUApiKey | <placeholder> | The Document Understanding API Key
#   6 OCREndpoint              | The OCR Engine endpoint

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 5 (shifts the old "Endpoint" row to row 6)
# so we can introduce the "DUApiKey" entry between "ApiKey" and "Endpoint".
$ws.Rows("5:5").Insert()

# Rename the existing "ApiKey" row to "OCRApiKey"
$ws.Range("A4").Value = "OCRApiKey"

# Fill in the newly inserted row with the "DUApiKey" entry
$ws.Range("A5").Value = "DUApiKey"

# Rename the old "Endpoint" row (now shifted to row 6) to "OCREndpoint"
$ws.Range("A6").Value = "OCREndpoint"

# Add the new "Description" column (C) with a bold header matching A1/B1
$ws.Range("C1").Value = "Description"
$ws.Range("C1").Font.Bold = $true

$ws.Range("C2").Value = "The path of the directory in which the files to be processed are stored"
$ws.Range("C3").Value = "The message that should be displayed when the workflow execution ends"
$ws.Range("C4").Value = "The OCR Engine API Key"
$ws.Range("C5").Value = "The Document Understanding API Key"
$ws.Range("C6").Value = "The OCR Engine endpoint"

# Widen the new Description column
$ws.Columns("C").ColumnWidth = 67.9

# B5 is left as an (empty) placeholder cell styled like a hyperlink
$ws.Range("B5").Style = "Hyperlink"
Write-Host "Applied Hyperlink style to B5"

# Match the saved selection state (A5 active cell)
[void]$ws.Range("A5").Select()
